# edit.ps1 - applies the "Time -> Biology" content rewrite described by the diff.

$d = $word.ActiveDocument

# Find-and-replace that writes the replacement text directly into the Range
# (rather than going through Find.Execute's replacement argument) so that
# Word's "smart quotes" autocorrect does not mangle apostrophes.
function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        Write-Host "WARNING: replace failed for: $old"
        return
    }
    $rng.Text = $new
}

# Locates $old and returns the Start position of the match, after clearing
# its text (used as an insertion anchor for building new run structure).
function Find-AndClear($old) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        Write-Host "WARNING: find failed for: $old"
        return -1
    }
    $pos = $rng.Start
    $rng.Text = ""
    return $pos
}

# Inserts $text (PowerShell backtick-v = soft line break) as a brand new run
# at $pos with the given Aptos/size/black formatting, returning the position
# right after the inserted text. Uses InsertAfter (not a Text= assignment)
# so that it always creates a new run boundary instead of merging into an
# adjacent run with identical formatting.
function Insert-Run($pos, $text, $fontSize) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    $newPos = $pos + $text.Length
    $fmtRange = $d.Range($pos, $newPos)
    $fmtRange.Font.Name = "Aptos"
    $fmtRange.Font.Size = $fontSize
    $fmtRange.Font.Color = 0
    return $newPos
}

# ---------------------------------------------------------------------------
# Title / byline / contact line
# ---------------------------------------------------------------------------
Replace-Text "Unveiling the Enigma of Time" "Exploring the Symphony of Life: A Journey through Biology"
Replace-Text "Dr. Jacob Templeton" "Anna Willis"
Replace-Text "jacob.templeton@academia.org" "at"

# ---------------------------------------------------------------------------
# Body paragraph (paragraph 5) - sentence by sentence rewrite
# ---------------------------------------------------------------------------
Replace-Text "Time, an omnipresent dimension, governs our perception of existence" "In an intricate dance of molecules and cells, the wonders of life unfold, weaving a symphony of wonder, knowledge and mystery"
Replace-Text " Its enigmatic nature has intrigued philosophers, physicists, and artists alike" " Biology, the science that unravels this enchanting world, invites us on a captivating journey to explore the very essence of living organisms"
Replace-Text " What is the essence of time? How do we measure its relentless passage? Can we manipulate or transcend its boundaries? These profound questions have spawned a wealth of theories, experiments, and artistic expressions, shaping our understanding of the universe and our place within it" " From the smallest building blocks of life to the diverse tapestry of ecosystems, this remarkable science reveals the beauty, complexity and interconnectedness of all living things"

Replace-Text "Early civilizations conceived time as a cyclical rhythm, mirroring the cycles of nature" "In this exploration, we'll venture into the microscopic realm of cells, discovering their intricate structures and remarkable functions"
Replace-Text " The sun's daily journey, the waxing and waning of the moon, and the seasonal shifts influenced ancient calendars and cosmologies" " We will unravel the genetic code, the blueprint of life, and witness the awe-inspiring process of DNA replication and protein synthesis"
Replace-Text " These early notions laid the foundation for the development of more sophisticated timekeeping methods, such as sundials, water clocks, and mechanical clocks, which gradually refined our perception of time's passage" " The dance of chromosomes during cell division, the intricate web of chemical reactions in metabolism, and the fascinating world of genetics will all come alive before our eyes"

Replace-Text "As science advanced, so did our understanding of time" "Moreover, we will delve into the depths of evolution, tracing the remarkable journey of life from its humble origins to the incredible diversity we see today"
Replace-Text " The 17th-century physicist Isaac Newton viewed time as an absolute and universal entity, ticking uniformly and independent of any observer" " We'll explore the mechanisms of natural selection, adaptation and speciation, understanding how life has adapted and thrived in the face of constant change"
Replace-Text " This Newtonian conception of time remained dominant for over two centuries, until Albert Einstein's groundbreaking theory of relativity shattered our conventional notions" " From the vast oceans to the lush forests, from towering mountains to arid deserts, we'll uncover the incredible adaptations that have allowed organisms to flourish in every corner of our planet"

# The final sentence of paragraph 5 ("Einstein revealed...") is removed and
# replaced with a large block of new runs/paragraph content. We locate and
# clear it, then rebuild the run structure manually so that line breaks and
# sentence-ending "." runs land exactly as in the target document. The
# original trailing "." run (already present after this sentence) is left
# untouched, so it will naturally follow whatever we insert here.
$pos = Find-AndClear " Einstein revealed that time is relative, intricately intertwined with space and affected by the presence of mass and energy"

$pos = Insert-Run $pos "`v" 12
$pos = Insert-Run $pos "`vIntroduction Continued:" 12
$pos = Insert-Run $pos "`v" 12
$pos = Insert-Run $pos "`vBiology extends its reach into the realm of ecosystems, where intricate webs of interdependence connect all living organisms" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos " We'll examine how species interact, forming dynamic communities that shape and sustain the delicate balance of life" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos " From predator-prey relationships to symbiotic partnerships, we'll unravel the intricate web of interactions that orchestrate the symphony of life" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos "`v" 12
$pos = Insert-Run $pos "`vFurthermore, we'll investigate the human body, a marvel of engineering, resilience and complexity" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos " We'll explore the intricacies of our organs and tissues, marveling at the intricate systems that regulate our heartbeat, respiration, digestion, and more" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos " The human body, with its intricate network of cells, tissues, and organs, provides a testament to the wonders of life's design" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos "`v" 12
$pos = Insert-Run $pos "`vFinally, we'll touch upon the frontiers of biology, where cutting-edge research is pushing the boundaries of our knowledge" 12
$pos = Insert-Run $pos "." 12
$pos = Insert-Run $pos " From advancements in genomics and biotechnology to the quest for understanding the origins of life, we'll glimpse the exciting possibilities that lie ahead" 12

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------
Replace-Text "Throughout history, time has been a central enigma, inspiring countless scientific investigations and artistic expressions" "Throughout this journey through biology, we'll unravel the mysteries of life, uncovering the elegance of cellular processes, the marvel of evolution, the interconnectedness of ecosystems and the intricacies of the human body"
Replace-Text " From ancient civilizations' cyclical perception of time to Newton's absolute time and Einstein's revolutionary theory of relativity, our understanding of time has undergone profound transformations" " We'll explore the frontiers of research, peering into the future of this captivating science"
Replace-Text " Time's elusive nature continues to challenge our intellect, fueling our quest for deeper knowledge and a comprehensive understanding of the fabric of reality" " Biology, a symphony of life, stands as a testament to the wonder, beauty and complexity of our universe"

# ---------------------------------------------------------------------------
# Trailing empty paragraph added at the very end of the body.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Host "Edit complete."
